$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row at 12, pushing the old totals row (row 12 -> row 13)
$ws1.Rows.Item(12).Insert()

# Existing row 11 client changes from "VIEJO RIVAS MAYRA ANABELLE" to
# "MOROCHO BACUILIMA HILDA INES" (advisor in column A is unchanged)
$ws1.Range("B11").Value = "MOROCHO BACUILIMA HILDA INES"

# New row 12 holds the client that used to be on row 11
$ws1.Range("A12").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B12").Value = "VIEJO RIVAS MAYRA ANABELLE"
$ws1.Range("C12:R12").Value = 0

# The totals row (now row 13) counts one more record: "0 de 10" -> "0 de 11", etc.
$ws1.Range("C13").Value = "0 de 11"
$ws1.Range("D13").Value = "0 de 11"
$ws1.Range("E13").Value = "0 de 11"
$ws1.Range("F13").Value = "0 de 11"
$ws1.Range("G13").Value = "0 de 11"
$ws1.Range("H13").Value = "0 de 11"
$ws1.Range("I13").Value = "0 de 11"
$ws1.Range("J13").Value = "0 de 11"
$ws1.Range("K13").Value = "0 de 11"
$ws1.Range("L13").Value = "0 de 11"
$ws1.Range("M13").Value = "3 de 11"
$ws1.Range("N13").Value = "0 de 11"
$ws1.Range("O13").Value = "0 de 11"
$ws1.Range("P13").Value = "1 de 11"
$ws1.Range("Q13").Value = "0 de 11"
$ws1.Range("R13").Value = "0 de 11"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item(2)

# Insert a new row at 12, pushing the old totals row (row 12 -> row 13)
$ws2.Rows.Item(12).Insert()

# Existing row 11 client changes from "VIEJO RIVAS MAYRA ANABELLE" to
# "MOROCHO BACUILIMA HILDA INES" (advisor in column A is unchanged)
$ws2.Range("B11").Value = "MOROCHO BACUILIMA HILDA INES"

# New row 12 holds the client that used to be on row 11
$ws2.Range("A12").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B12").Value = "VIEJO RIVAS MAYRA ANABELLE"
$ws2.Range("C12:G12").Value = 0

# The totals row (now row 13) keeps the same numeric totals as before
$ws2.Range("C13").Value = 1187.62
$ws2.Range("D13").Value = 0
$ws2.Range("E13").Value = 0
$ws2.Range("F13").Value = 5551.6
$ws2.Range("G13").Value = 1200
